$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (column D becomes "Origin Airport", H "Sold seats" column is dropped later)
$ws.Range("A1").Value = "Airlline code"
$ws.Range("B1").Value = "Flight nbr"
$ws.Range("C1").Value = "Dept date"
$ws.Range("D1").Value = "Origin Airport"
$ws.Range("E1").Value = "Destination Airport"
$ws.Range("F1").Value = "Cabin "
$ws.Range("G1").Value = "Empty seats"

# Data rows: Airline, Flight, Dept date(serial), Origin, Destination, Cabin, Empty seats
$data = @(
    @("SQ", "SQ1234", 43639, "KUL", "COK", "Y", 10),
    @("SQ", "SQ1234", 43639, "KUL", "COK", "W", 23),
    @("SQ", "SQ1234", 43639, "KUL", "COK", "C", 22),
    @("SQ", "SQ1235", 43639, "SIN", "SFO", "Y", 34),
    @("SQ", "SQ1235", 43639, "SIN", "SFO", "W", 33),
    @("SQ", "SQ1235", 43639, "SIN", "SFO", "C", 45),
    @("SQ", "SQ1236", 43639, "HND", "SFO", "Y", 33),
    @("SQ", "SQ1236", 43639, "HND", "SFO", "W", 12),
    @("SQ", "SQ1236", 43639, "HND", "SFO", "C", 11)
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $r = $r + 1
}

# Remove the now-unused "Sold seats" column (H)
$ws.Columns.Item(8).Delete()

# Update the selected cell to match the saved workbook state
$ws.Range("E9").Select() | Out-Null
